$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2941544.8
$ws.Range("I9").Value = 4902374
$ws.Range("K9").Value = 4902374
$ws.Range("M9").Value = -4902205
$ws.Range("H15").Value = 2601.8416
$ws.Range("I15").Value = 2601.8416
$ws.Range("K15").Value = 7805.524800000001
$ws.Range("M15").Value = -7636.524800000001
$ws.Range("H40").Value = 4583.9546
$ws.Range("I40").Value = 5120
$ws.Range("J40").Value = 4212.846
$ws.Range("K40").Value = 5120
$ws.Range("L40").Value = 4212.846
$ws.Range("M40").Value = -4945
$ws.Range("N40").Value = -4562.846
$ws.Range("H70").Value = 784
$ws.Range("J70").Value = 925
$ws.Range("L70").Value = 2775
$ws.Range("N70").Value = -3315
$ws.Range("H73").Value = 784
$ws.Range("J73").Value = 925
$ws.Range("L73").Value = 2775
$ws.Range("N73").Value = -4647
$ws.Range("H86").Value = 222228940
$ws.Range("I86").Value = 333341900
$ws.Range("K86").Value = 333341900
$ws.Range("M86").Value = -333340777
$ws.Range("H89").Value = 222228940
$ws.Range("I89").Value = 333341900
$ws.Range("K89").Value = 1666709500
$ws.Range("M89").Value = -1666703884
$ws.Range("H107").Value = 807.8333
$ws.Range("J107").Value = 400
$ws.Range("L107").Value = 400
$ws.Range("N107").Value = -4240
$ws.Range("H132").Value = 10340.375
$ws.Range("I132").Value = 10340.375
$ws.Range("K132").Value = 31021.125
$ws.Range("M132").Value = -28491.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 565.93335
$ws.Range("I2").Value = 549.7083
$ws.Range("K2").Value = 549.7083
$ws.Range("M2").Value = -436.7083
$ws.Range("H32").Value = 5727.467
$ws.Range("I32").Value = 4353.125
$ws.Range("J32").Value = 11224.833
$ws.Range("K32").Value = 4353.125
$ws.Range("L32").Value = 11224.833
$ws.Range("M32").Value = -4066.125
$ws.Range("N32").Value = -11798.833
$ws.Range("H61").Value = 11657.223
$ws.Range("I61").Value = 3500
$ws.Range("K61").Value = 3500
$ws.Range("M61").Value = -3288
$ws.Range("H116").Value = 565.93335
$ws.Range("I116").Value = 549.7083
$ws.Range("K116").Value = 549.7083
$ws.Range("M116").Value = 1744.2917
$ws.Range("H136").Value = 11657.223
$ws.Range("I136").Value = 3500
$ws.Range("K136").Value = 10500
$ws.Range("M136").Value = -7950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 565.93335
$ws.Range("I3").Value = 549.7083
$ws.Range("K3").Value = 549.7083
$ws.Range("M3").Value = -435.7083
$ws.Range("H20").Value = 26886192
$ws.Range("I20").Value = 32055986
$ws.Range("K20").Value = 32055986
$ws.Range("M20").Value = -32055739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4233.4565
$ws.Range("I31").Value = 3716.1365
$ws.Range("J31").Value = 4707.6665
$ws.Range("K31").Value = 3716.1365
$ws.Range("L31").Value = 4707.6665
$ws.Range("M31").Value = -3421.1365
$ws.Range("N31").Value = -5297.6665
$ws.Range("H34").Value = 4233.4565
$ws.Range("I34").Value = 3716.1365
$ws.Range("J34").Value = 4707.6665
$ws.Range("K34").Value = 3716.1365
$ws.Range("L34").Value = 4707.6665
$ws.Range("M34").Value = -3514.1365
$ws.Range("N34").Value = -5111.6665
$ws.Range("H58").Value = 4332
$ws.Range("I58").Value = 4524.5
$ws.Range("J58").Value = 4272.769
$ws.Range("K58").Value = 4524.5
$ws.Range("L58").Value = 4272.769
$ws.Range("M58").Value = -4321.5
$ws.Range("N58").Value = -4678.769
$ws.Range("H63").Value = 95000
$ws.Range("J63").Value = 95000
$ws.Range("L63").Value = 95000
$ws.Range("N63").Value = -96372
$ws.Range("H66").Value = 95000
$ws.Range("J66").Value = 95000
$ws.Range("L66").Value = 285000
$ws.Range("N66").Value = -291864
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H107").Value = 4167710.8
$ws.Range("I107").Value = 4546548
$ws.Range("K107").Value = 4546548
$ws.Range("M107").Value = -4544628
$ws.Range("H122").Value = 4363.9287
$ws.Range("I122").Value = 4154.1
$ws.Range("K122").Value = 12462.3
$ws.Range("M122").Value = -10012.3
$ws.Range("H136").Value = 4332
$ws.Range("I136").Value = 4524.5
$ws.Range("J136").Value = 4272.769
$ws.Range("K136").Value = 13573.5
$ws.Range("L136").Value = 12818.307
$ws.Range("M136").Value = -11023.5
$ws.Range("N136").Value = -17918.307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 546.2273
$ws.Range("J107").Value = 624.375
$ws.Range("L107").Value = 1873.125
$ws.Range("N107").Value = -5713.125
$ws.Range("H132").Value = 2704.425
$ws.Range("I132").Value = 1683.25
$ws.Range("J132").Value = 3142.0715
$ws.Range("K132").Value = 15149.25
$ws.Range("L132").Value = 28278.6435
$ws.Range("M132").Value = -12619.25
$ws.Range("N132").Value = -33338.6435
$ws.Range("H134").Value = 3632.6667
$ws.Range("I134").Value = 3632.6667
$ws.Range("K134").Value = 10898.0001
$ws.Range("M134").Value = -5828.000100000001
$ws.Range("H140").Value = 9902.968999999999
$ws.Range("I140").Value = 4813.4546
$ws.Range("K140").Value = 14440.3638
$ws.Range("M140").Value = -9260.363799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13572112
$ws.Range("I70").Value = 20920820
$ws.Range("J70").Value = 5265.923
$ws.Range("K70").Value = 20920820
$ws.Range("L70").Value = 5265.923
$ws.Range("M70").Value = -20920550
$ws.Range("N70").Value = -5805.923
$ws.Range("H73").Value = 13572112
$ws.Range("I73").Value = 20920820
$ws.Range("J73").Value = 5265.923
$ws.Range("K73").Value = 20920820
$ws.Range("L73").Value = 5265.923
$ws.Range("M73").Value = -20919884
$ws.Range("N73").Value = -7137.923
$ws.Range("H102").Value = 1237.5
$ws.Range("J102").Value = 2097.625
$ws.Range("L102").Value = 2097.625
$ws.Range("N102").Value = -5341.625
$ws.Range("H107").Value = 996.6667
$ws.Range("I107").Value = 996.6667
$ws.Range("K107").Value = 996.6667
$ws.Range("M107").Value = 923.3333
$ws.Range("H133").Value = 134000
$ws.Range("J133").Value = 134000
$ws.Range("L133").Value = 134000
$ws.Range("N133").Value = -144120
$ws.Range("H135").Value = 68873.625
$ws.Range("J135").Value = 68873.625
$ws.Range("L135").Value = 68873.625
$ws.Range("N135").Value = -79013.625
$ws.Range("H136").Value = 10336.444
$ws.Range("J136").Value = 10336.444
$ws.Range("L136").Value = 31009.332
$ws.Range("N136").Value = -36109.33199999999
$ws.Range("H140").Value = 73384.94
$ws.Range("J140").Value = 73384.94
$ws.Range("L140").Value = 73384.94
$ws.Range("N140").Value = -83744.94

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 542.3077
$ws.Range("I22").Value = 510.7143
$ws.Range("J22").Value = 579.1667
$ws.Range("K22").Value = 510.7143
$ws.Range("L22").Value = 579.1667
$ws.Range("M22").Value = -215.7143
$ws.Range("N22").Value = -1169.1667
$ws.Range("H27").Value = 542.3077
$ws.Range("I27").Value = 510.7143
$ws.Range("J27").Value = 579.1667
$ws.Range("K27").Value = 510.7143
$ws.Range("L27").Value = 579.1667
$ws.Range("M27").Value = -403.7143
$ws.Range("N27").Value = -793.1667
$ws.Range("H93").Value = 1492.3158
$ws.Range("I93").Value = 1492.3158
$ws.Range("K93").Value = 1492.3158
$ws.Range("M93").Value = -244.3158000000001
$ws.Range("H122").Value = 3991.3333
$ws.Range("I122").Value = 3975
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 11925
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -9475
$ws.Range("N122").Value = -16898.5
$ws.Range("H136").Value = 5178.2964
$ws.Range("J136").Value = 5091.0835
$ws.Range("L136").Value = 15273.2505
$ws.Range("N136").Value = -20373.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 12824285
$ws.Range("I132").Value = 17548136
$ws.Range("J132").Value = 2407
$ws.Range("K132").Value = 52644408
$ws.Range("L132").Value = 7221
$ws.Range("M132").Value = -52641878
$ws.Range("N132").Value = -12281
